# Update column G ("K") values on Sheet1 to reflect the regenerated
# strike-count data (commit: "regen save_data to use K instead of Strike#,
# regen std/mean, calc and write s_vals").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 1
    3  = 2
    4  = 2
    5  = 3
    6  = 0
    7  = 1
    8  = 1
    9  = 2
    10 = 2
    11 = 4
    12 = 3
    13 = 2
    14 = 3
    15 = 3
    16 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
